# Update the workbook according to the target diff:
# 1. Column C (Förändrad) for rows 2-438: 45172 -> 45175
# 2. Row 438 gets an explicit row height (15pt, custom height)
# 3. A new row 439 is appended with a new entry
# 4. Sheet dimension grows to A1:Y439 (handled automatically by Excel)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Bulk-update column C for all existing data rows (2 through 438)
$ws.Range("C2:C438").Value = 45175

# 2. Make sure row 438 has the same explicit 15pt row height as the rows above it
$ws.Rows.Item(438).RowHeight = 15

# 3. Append the new row of data (row 439)
$ws.Cells.Item(439, 1).Value = "A 41182-2023"

$ws.Cells.Item(439, 2).Value = 45174
$ws.Cells.Item(439, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(439, 3).Value = 45175
$ws.Cells.Item(439, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(439, 4).Value = "UPPSALA LÄN"
$ws.Cells.Item(439, 5).Value = "TIERP"
$ws.Cells.Item(439, 6).Value = "Bergvik skog öst AB"

$ws.Cells.Item(439, 7).Value = 14.8
$ws.Cells.Item(439, 8).Value = 0
$ws.Cells.Item(439, 9).Value = 0
$ws.Cells.Item(439, 10).Value = 0
$ws.Cells.Item(439, 11).Value = 0
$ws.Cells.Item(439, 12).Value = 0
$ws.Cells.Item(439, 13).Value = 0
$ws.Cells.Item(439, 14).Value = 0
$ws.Cells.Item(439, 15).Value = 0
$ws.Cells.Item(439, 16).Value = 0
$ws.Cells.Item(439, 17).Value = 0

# Column R (18) keeps the same "wrap text" style used throughout the sheet,
# with no content, matching the other rows.
$ws.Cells.Item(439, 18).WrapText = $true

Write-Host "Edit complete"
Write-Host ("C2=" + $ws.Range("C2").Value())
Write-Host ("C438=" + $ws.Range("C438").Value())
Write-Host ("A439=" + $ws.Cells.Item(439,1).Value())
Write-Host ("B439=" + $ws.Cells.Item(439,2).Value())
Write-Host ("C439=" + $ws.Cells.Item(439,3).Value())
Write-Host ("G439=" + $ws.Cells.Item(439,7).Value())
